$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, matching the style used by the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Compute "Save" flag for each data row: 1 if the "sum" (column G) exceeds 8, else 0
for ($row = 2; $row -le 50; $row++) {
    $sumVal = $ws.Cells.Item($row, 7).Value()
    if ($sumVal -gt 8) {
        $ws.Cells.Item($row, 8).Value = 1
    } else {
        $ws.Cells.Item($row, 8).Value = 0
    }
}
